$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.936.00"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.637.80"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.96%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "1.864.94"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").Value = "1.578.21"
$ws.Range("E14").Value = "  -4.74%  "
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "25.956.15"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("E27").Value = "  +2.77%  "
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.901"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").Value = "1.136.13"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.544"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.797"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("D45").Value = "1.774.04"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("E46").Value = "  +5.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.415"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.37%  "
